# Realestate Update resale numbers 2024-01-03 17:20
# Appends a new observation row (row 12) to the CityResaleNum sheet with
# the resale-index readings captured at 2024-01-03 17:20:14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Date / "00" week columns look like a date and a number to Excel's normal
# typed-input parser, so mark them as Text first to preserve the literal
# strings (matches how the rest of the sheet stores them).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-03"

$ws.Cells.Item($row, 2).Value = "17:20:14"
$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"

$ws.Cells.Item($row, 5).Value = 140047
$ws.Cells.Item($row, 6).Value = 142802
$ws.Cells.Item($row, 7).Value = 172075
$ws.Cells.Item($row, 8).Value = 146914
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 117375
$ws.Cells.Item($row, 11).Value = 223977
$ws.Cells.Item($row, 12).Value = 248039
$ws.Cells.Item($row, 13).Value = 183895
$ws.Cells.Item($row, 14).Value = 109860
$ws.Cells.Item($row, 15).Value = 40107
$ws.Cells.Item($row, 16).Value = 30823
$ws.Cells.Item($row, 17).Value = 72174
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41489
$ws.Cells.Item($row, 20).Value = -1
